# Update "想去人数" (number of people interested) values (column F)
# across the "展览", "演出" and "全部类型" sheets, reflecting refreshed
# output data (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 809
$ws1.Range("F4").Value = 1109
$ws1.Range("F5").Value = 177
$ws1.Range("F7").Value = 222
$ws1.Range("F8").Value = 394
$ws1.Range("F9").Value = 1015
$ws1.Range("F10").Value = 13
$ws1.Range("F11").Value = 510
$ws1.Range("F13").Value = 159
$ws1.Range("F14").Value = 12696
$ws1.Range("F15").Value = 5
$ws1.Range("F16").Value = 5232

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 89

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 809
$ws4.Range("F5").Value = 1109
$ws4.Range("F6").Value = 177
$ws4.Range("F8").Value = 222
$ws4.Range("F9").Value = 394
$ws4.Range("F10").Value = 1015
$ws4.Range("F11").Value = 13
$ws4.Range("F12").Value = 510
$ws4.Range("F14").Value = 159
$ws4.Range("F15").Value = 12696
$ws4.Range("F16").Value = 89
$ws4.Range("F18").Value = 5
$ws4.Range("F19").Value = 5232
